$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Data")

# Extend the data table by one row (row 19), carrying over row 18's
# cell formatting (date/number styles), and fill in the new values.
$ws.Range("A18:O18").Copy()
$ws.Range("A19:O19").PasteSpecial(-4122) # xlPasteFormats
$excel.CutCopyMode = $false

$ws.Cells.Item(19, 1).Value = 45728.502708333333  # TimeStamp
$ws.Cells.Item(19, 2).Value = 10
$ws.Cells.Item(19, 3).Value = 6
$ws.Cells.Item(19, 4).Value = 240
$ws.Cells.Item(19, 5).Value = 426
$ws.Cells.Item(19, 6).Value = 402
$ws.Cells.Item(19, 7).Value = 476
$ws.Cells.Item(19, 8).Value = 3432
$ws.Cells.Item(19, 9).Value = 476
$ws.Cells.Item(19, 10).Value = 2026
$ws.Cells.Item(19, 11).Value = 208
$ws.Cells.Item(19, 12).Value = 417
$ws.Cells.Item(19, 13).Value = 30
$ws.Cells.Item(19, 14).Value = 3683
$ws.Cells.Item(19, 15).Value = 4725

$ws.Range("E18").Select()
